$d = $word.ActiveDocument

$find = $d.Content.Find
$find.Execute("plaintiff_name", $true, $false, $false, $false, $false, $true, 1, $false, "opposing_party.name", 2)
